$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.148.79'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.24%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.658.13'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.14%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.11'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5192'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.35%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.33%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2627'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -2.64%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06266'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.76'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -4.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07711'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.15%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.426'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.96%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.645.51'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.87%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.884.19'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5414'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.20%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8128'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.64%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.67'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.58%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.177.40'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.21%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.39%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.619'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.61%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '191.37'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.06'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.35%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.035'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -4.68%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.45%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '139.27'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.62%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1227'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.176'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -3.22%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.05'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.402'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -3.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05956'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.268'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.82%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.544'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.93%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.253'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -5.92%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.603'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -5.53%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9648'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -4.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.427'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.28%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.771'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5678'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -8.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01592'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.51%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.968'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.18%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8552'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.82%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.38%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.009.30'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -7.99%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.39'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.799.38'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.17%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₈109'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '56.58'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.67%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.007'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.18%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.978'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05174'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.40%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4197'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.86%  '
